# R301_Grille_Evaluation.xlsx - "mise a jour grille eval"
# Update the "nous" column (D) self-evaluation scores on Feuille1:
# several sub-items go from 0 -> 0.5. The SUM() subtotal/total formulas
# in column D (D13, D20, D26, D47) recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuille1")

$ws.Range("D8").Value  = 0.5
$ws.Range("D10").Value = 0.5
$ws.Range("D11").Value = 0.5
$ws.Range("D12").Value = 0.5
$ws.Range("D17").Value = 0.5
$ws.Range("D18").Value = 0.5

# Update the view state: scroll so row 2 is the top-left visible row,
# zoom out to 65%, and move the active selection to R11.
$ws.Activate()
$excel.ActiveWindow.DisplayGridlines = $true
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 65
$ws.Range("R11").Select() | Out-Null
